{"js": "// Remove the \" (m/w)\" suffix from the job title heading\n// (\"Director of Mid Market Accounts (m/w)\" -> \"Director of Mid Market Accounts\")\nconst body = context.document.body;\nconst results = body.search(\" Accounts (m/w)\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\" Accounts\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Remove the \" (m/w)\" suffix from the job title heading\n# (\"Director of Mid Market Accounts (m/w)\" -> \"Director of Mid Market Accounts\")\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \" Accounts (m/w)\"\n$find.Replacement.Text = \" Accounts\"\n$find.Forward = $true\n$find.Wrap = 1          # wdFindContinue\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n# wdReplaceAll = 2\n$find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $false, $false, $find.Forward, $find.Wrap, $false, $find.Replacement.Text, 2) | Out-Null\n"}
